# Daily attendance processing - 2025-10-09 16:51:39
# Reverses the order of comma-separated "Recorded By" names in column G
# for rows that contain multiple recorders.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    2   = "backup@backdoor.com, System, system"
    3   = "dnasr281@gmail.com, System"
    4   = "backup@backdoor.com, System"
    5   = "backup@backdoor.com, System"
    6   = "dnasr281@gmail.com, System"
    11  = "dnasr281@gmail.com, System"
    12  = "dnasr281@gmail.com, System"
    13  = "dnasr281@gmail.com, System"
    14  = "dnasr281@gmail.com, System"
    15  = "dnasr281@gmail.com, System"
    29  = "backup@backdoor.com, System, system"
    30  = "dnasr281@gmail.com, System"
    32  = "backup@backdoor.com, System"
    33  = "dnasr281@gmail.com, System"
    38  = "dnasr281@gmail.com, System"
    39  = "dnasr281@gmail.com, System"
    40  = "dnasr281@gmail.com, System"
    41  = "dnasr281@gmail.com, System"
    42  = "dnasr281@gmail.com, System"
    56  = "backup@backdoor.com, System, system"
    57  = "dnasr281@gmail.com, System"
    58  = "backup@backdoor.com, System"
    59  = "backup@backdoor.com, System"
    60  = "dnasr281@gmail.com, System"
    65  = "dnasr281@gmail.com, System"
    66  = "dnasr281@gmail.com, System"
    67  = "dnasr281@gmail.com, System"
    68  = "dnasr281@gmail.com, System"
    69  = "dnasr281@gmail.com, System"
    84  = "backup@backdoor.com, System"
    85  = "backup@backdoor.com, System"
    86  = "dnasr281@gmail.com, System"
    89  = "dnasr281@gmail.com, System"
    90  = "dnasr281@gmail.com, admin@admin.com"
    93  = "dnasr281@gmail.com, System"
    95  = "dnasr281@gmail.com, System"
    110 = "backup@backdoor.com, System"
    111 = "backup@backdoor.com, System"
    112 = "dnasr281@gmail.com, System"
    115 = "dnasr281@gmail.com, System"
    116 = "dnasr281@gmail.com, admin@admin.com"
    119 = "dnasr281@gmail.com, System"
    121 = "dnasr281@gmail.com, System"
    136 = "backup@backdoor.com, System"
    137 = "backup@backdoor.com, System"
    138 = "dnasr281@gmail.com, System"
    141 = "dnasr281@gmail.com, System"
    142 = "dnasr281@gmail.com, admin@admin.com"
    145 = "dnasr281@gmail.com, System"
    147 = "dnasr281@gmail.com, System"
}

foreach ($rowNum in $updates.Keys) {
    $ws.Cells.Item($rowNum, 7).Value = $updates[$rowNum]
}
